# Apply the "2.0.0" term update to ValueSet-KLTheraphyAreasFSIII.xlsx
#
# Summary of changes (per commit "adding term 2.0.0 with four more concepts in FBOE"):
#  - Sheet "Metadata": bump Version, Date and Contact property values.
#  - Sheet "Include from FSIII": insert a new concept row (new UUID,
#    blank display column) above the trailing blank-row / "System URI"
#    footer rows, shifting those two rows down by one.

$wb = $excel.ActiveWorkbook

# --- Sheet: Metadata -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")

$ws1.Range("B3").Value  = "2.0.0"                                          # Version
$ws1.Range("B8").Value  = "2024-06-04T14:59:10+02:00"                      # Date
$ws1.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"       # Contact

# --- Sheet: Include from FSIII ---------------------------------------
$ws2 = $wb.Worksheets.Item("Include from FSIII")

# Make room for the new concept row: shift the last two existing rows
# (the blank-concept row and the "System URI" footer row) down by one,
# row 14 -> row 15 first (so we don't overwrite data we still need).
$ws2.Range("A12:B12").Copy()
$ws2.Range("A15:B15").PasteSpecial(-4122)
$ws2.Range("A15").Value = $ws2.Range("A14").Value()
$ws2.Range("B15").Value = $ws2.Range("B14").Value()

# row 13 -> row 14
$ws2.Range("A14").Value = $ws2.Range("A13").Value()
$ws2.Range("B14").Value = ""

# Write the new concept into the freed-up row 13.
$ws2.Range("A13").Value = "aec684bd-c2ea-4ff0-8eb7-6d2cf67fb863"
$ws2.Range("B13").ClearContents()
